$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated RelativeTime (B), AbsoluteTime... (C), Difference (D) values for rows 2-16
$data = @(
    @(2,  1.9995574000058696, 83803.127819800007, 83805.129512700005),
    @(3,  3.998875000004773,  83803.128058300004, 83807.128830300004),
    @(4,  5.9985075000004144, 83803.128059099996, 83809.128462799999),
    @(5,  7.9992842000065139, 83803.128056300004, 83811.129239500005),
    @(6,  9.9985098000033759, 83803.128059199997, 83813.128465100002),
    @(7,  11.998405100006494, 83803.128059200011, 83815.128360400005),
    @(8,  13.998957000003429, 83803.128057199996, 83817.128912300002),
    @(9,  15.998801499998081, 83803.128057099995, 83819.128756799997),
    @(10, 17.998547299997881, 83803.128059199997, 83821.128502599997),
    @(11, 19.998367299995152, 83803.128059499999, 83823.128322599994),
    @(12, 21.998547600000165, 83803.128058300004, 83825.128502899999),
    @(13, 23.998453299995163, 83803.128058899994, 83827.128408599994),
    @(14, 25.998485300005996, 83803.128058900009, 83829.128440600005),
    @(15, 27.998714700006531, 83803.128057599999, 83831.128670000006),
    @(16, 29.998485899996012, 83803.128059099996, 83833.128441199995)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Rows 17, 18, 19: clear out the old event data - A/B/C/D become 0, E becomes empty
foreach ($r in 17..19) {
    $ws.Cells.Item($r, 1).Value = 0
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = ""
}
